# Sudoku augmented example update:
#  - "checks" sheet gains a third column ("pots") with per-check pot counts
#  - "entries" sheet keeps its existing "pots" column (unchanged values)
#  - the active sheet/selection moves from "entries" back to "checks"

$wb = $excel.ActiveWorkbook
$checks = $wb.Worksheets.Item("checks")
$entries = $wb.Worksheets.Item("entries")

# --- checks: add the "pots" column ------------------------------------
# Header style matches the existing header cells (A1/B1).
$checks.Range("B1").Copy()
$checks.Range("C1").PasteSpecial(-4122)  # xlPasteFormats

# Numeric-cell style matches the "pots" column already used on "entries".
$entries.Range("C2").Copy()
$checks.Range("C2:C4").PasteSpecial(-4122)  # xlPasteFormats

$checks.Range("C1").Value = "pots"
$checks.Range("C2").Value = 9
$checks.Range("C3").Value = 7
$checks.Range("C4").Value = 2

# --- selections / active sheet -----------------------------------------
$entries.Range("C1").Select() | Out-Null
$checks.Range("C4").Select() | Out-Null

$checks.Activate()
